{"js": "const tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nconst edits = [\n  { row: 0, col: 0, oldText: \"24\u00f74=\", newText: \"66\u00f78=\" },\n  { row: 0, col: 1, oldText: \"97\u00f73=\", newText: \"15\u00f77=\" },\n  { row: 0, col: 2, oldText: \"33\u00f78=\", newText: \"61\u00f77=\" },\n  { row: 0, col: 3, oldText: \"56\u00f79=\", newText: \"10\u00f76=\" },\n  { row: 0, col: 4, oldText: \"35\u00f76=\", newText: \"54\u00f75=\" },\n  { row: 4, col: 0, oldText: \"74\u00f73=\", newText: \"49\u00f79=\" },\n  { row: 4, col: 1, oldText: \"22\u00f75=\", newText: \"54\u00f73=\" },\n  { row: 4, col: 2, oldText: \"76\u00f78=\", newText: \"94\u00f73=\" },\n  { row: 4, col: 3, oldText: \"22\u00f72=\", newText: \"41\u00f79=\" },\n  { row: 4, col: 4, oldText: \"19\u00f75=\", newText: \"45\u00f74=\" },\n  { row: 8, col: 0, oldText: \"76\u00f74=\", newText: \"50\u00f77=\" },\n  { row: 8, col: 1, oldText: \"40\u00f73=\", newText: \"56\u00f79=\" },\n  { row: 8, col: 2, oldText: \"32\u00f79=\", newText: \"18\u00f74=\" },\n  { row: 8, col: 3, oldText: \"50\u00f78=\", newText: \"94\u00f78=\" },\n  { row: 8, col: 4, oldText: \"55\u00f73=\", newText: \"66\u00f73=\" },\n  { row: 12, col: 0, oldText: \"90\u00f79=\", newText: \"97\u00f78=\" },\n  { row: 12, col: 1, oldText: \"22\u00f74=\", newText: \"36\u00f77=\" },\n  { row: 12, col: 2, oldText: \"67\u00f76=\", newText: \"60\u00f73=\" },\n  { row: 12, col: 3, oldText: \"67\u00f73=\", newText: \"50\u00f77=\" },\n  { row: 12, col: 4, oldText: \"15\u00f79=\", newText: \"44\u00f77=\" },\n  { row: 16, col: 0, oldText: \"31\u00f75=\", newText: \"54\u00f79=\" },\n  { row: 16, col: 1, oldText: \"27\u00f77=\", newText: \"73\u00f77=\" },\n  { row: 16, col: 2, oldText: \"69\u00f77=\", newText: \"42\u00f79=\" },\n  { row: 16, col: 3, oldText: \"19\u00f72=\", newText: \"65\u00f76=\" },\n  { row: 16, col: 4, oldText: \"61\u00f77=\", newText: \"12\u00f74=\" },\n];\n\nconst paras = [];\nfor (const e of edits) {\n  const cell = table.getCell(e.row, e.col);\n  cell.body.paragraphs.load(\"items,text\");\n  paras.push(cell.body.paragraphs);\n}\nawait context.sync();\n\nfor (let i = 0; i < edits.length; i++) {\n  const para = paras[i].items[0];\n  // Replace the cell's paragraph text with the new problem. Cells are targeted\n  // by (row, col) position rather than by searching for the old text, since\n  // some \"new\" values equal \"old\" values of other cells (e.g. \"61\u00f77=\" and\n  // \"56\u00f79=\" both appear as an old value in one cell and a new value in\n  // another), which would make a sequential text search-and-replace unsafe.\n  para.insertText(edits[i].newText, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Update the 25 \"two-digit \u00f7 one-digit\" practice problems in place, addressed\n# by (row, col) so the new values (which sometimes equal another cell's old\n# value, e.g. \"61\u00f77=\") can't collide with a text-search based replace.\n$t.Cell(1, 1).Range.Text = \"66\u00f78=\"\n$t.Cell(1, 2).Range.Text = \"15\u00f77=\"\n$t.Cell(1, 3).Range.Text = \"61\u00f77=\"\n$t.Cell(1, 4).Range.Text = \"10\u00f76=\"\n$t.Cell(1, 5).Range.Text = \"54\u00f75=\"\n$t.Cell(5, 1).Range.Text = \"49\u00f79=\"\n$t.Cell(5, 2).Range.Text = \"54\u00f73=\"\n$t.Cell(5, 3).Range.Text = \"94\u00f73=\"\n$t.Cell(5, 4).Range.Text = \"41\u00f79=\"\n$t.Cell(5, 5).Range.Text = \"45\u00f74=\"\n$t.Cell(9, 1).Range.Text = \"50\u00f77=\"\n$t.Cell(9, 2).Range.Text = \"56\u00f79=\"\n$t.Cell(9, 3).Range.Text = \"18\u00f74=\"\n$t.Cell(9, 4).Range.Text = \"94\u00f78=\"\n$t.Cell(9, 5).Range.Text = \"66\u00f73=\"\n$t.Cell(13, 1).Range.Text = \"97\u00f78=\"\n$t.Cell(13, 2).Range.Text = \"36\u00f77=\"\n$t.Cell(13, 3).Range.Text = \"60\u00f73=\"\n$t.Cell(13, 4).Range.Text = \"50\u00f77=\"\n$t.Cell(13, 5).Range.Text = \"44\u00f77=\"\n$t.Cell(17, 1).Range.Text = \"54\u00f79=\"\n$t.Cell(17, 2).Range.Text = \"73\u00f77=\"\n$t.Cell(17, 3).Range.Text = \"42\u00f79=\"\n$t.Cell(17, 4).Range.Text = \"65\u00f76=\"\n$t.Cell(17, 5).Range.Text = \"12\u00f74=\"\n"}
